$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 451; existing rows 451..537 shift down to 452..538
$ws.Rows.Item(451).Insert()

# Populate the newly inserted row 451 with the new record
$ws.Range("A451").Value = 4
$ws.Range("B451").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C451").Value = "Los Lagos"
$ws.Range("D451").Value = 45211
$ws.Range("E451").Value = 10
$ws.Range("F451").Value = 100112040
$ws.Range("G451").Value = "Cilantro"
$ws.Range("H451").Value = "Sin especificar"
$ws.Range("I451").Value = "Primera"
$ws.Range("J451").Value = 100
$ws.Range("K451").Value = 12000
$ws.Range("L451").Value = 12000
$ws.Range("M451").Value = 12000
$ws.Range("N451").Value = '$/caja 36 atados'
$ws.Range("O451").Value = "Región Metropolitana"
$ws.Range("P451").Value = 333
$ws.Range("Q451").Value = 36
$ws.Range("R451").Value = "Hortaliza"
